# "Skip the header row" - insert a new header row at the top of the sheet
# with column titles (Name, Email, Age, Salary, Department), bold-formatted,
# pushing all existing data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remember the existing hyperlinks (cell -> mailto target) before we touch
#    anything; this simulator does not shift hyperlink ranges automatically
#    when rows are inserted, so we will recreate them afterwards.
$linkTargets = @(
    @{Row = 15; Target = "mailto:Fifteen@fiv.com"},
    @{Row = 14; Target = "mailto:Fourteen@one.com"},
    @{Row = 13; Target = "mailto:Thirteen@cu.com"},
    @{Row = 12; Target = "mailto:Abc123@arw.com"},
    @{Row = 10; Target = "mailto:Ten@c.com"},
    @{Row = 9;  Target = "mailto:Nine@c.com"},
    @{Row = 1;  Target = "mailto:N1@a.com"},
    @{Row = 3;  Target = "mailto:N3@"},
    @{Row = 5;  Target = "mailto:N5@.com"},
    @{Row = 8;  Target = "mailto:Eight@.com"},
    @{Row = 2;  Target = "mailto:Eight@.com"},
    @{Row = 4;  Target = "mailto:Eight@.com"},
    @{Row = 11; Target = "mailto:Eight@.com"},
    @{Row = 6;  Target = "mailto:name@.com"},
    @{Row = 7;  Target = "mailto:name@.com"}
)

# 2) Drop the old hyperlinks so we can rebuild them, shifted, afterwards.
$ws.Hyperlinks.Delete()

# 3) Insert a blank row above row 1; everything currently on the sheet
#    (rows 1-15) shifts down to rows 2-16.
$ws.Rows.Item(1).Insert()

# 4) Fill in the new header row with bold labels matching each column's
#    existing number format (text / integer / 2-decimal).
$ws.Range("A1").Value = "Name"
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Font.Bold = $true

$ws.Range("B1").Value = "Email"
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Font.Bold = $true

$ws.Range("C1").Value = "Age"
$ws.Range("C1").NumberFormat = "0"
$ws.Range("C1").Font.Bold = $true

$ws.Range("D1").Value = "Salary"
$ws.Range("D1").NumberFormat = "0.00"
$ws.Range("D1").Font.Bold = $true

$ws.Range("E1").Value = "Department"
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Font.Bold = $true

# 5) Re-create the hyperlinks on column B, one row lower than before.
#    Adding a hyperlink re-applies the "Hyperlink" look, which can nudge the
#    cell onto a freshly-minted style; restore the original text format /
#    centered alignment so the cell keeps using the existing shared style.
foreach ($link in $linkTargets) {
    $newRow = $link.Row + 1
    $ws.Hyperlinks.Add($ws.Range("B$newRow"), $link.Target) | Out-Null
    $ws.Range("B$newRow").NumberFormat = "@"
    $ws.Range("B$newRow").HorizontalAlignment = -4108
}

# 6) Select the new header row, matching the post-edit selection state.
$ws.Range("A1:E1").Select()
